$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.186.51"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "2.599.40"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.60"
$ws.Range("E5").Value = "  +4.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.08"
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.44"
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.102"
$ws.Range("E10").Value = "  +1.75%  "
$ws.Range("E11").Value = "  +1.20%  "
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("D13").Value = "3.063.22"
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("D14").Value = "59.121.46"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.52"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000133"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.587.95"
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "343.27"
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("E21").Value = "  -2.13%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.59"
$ws.Range("E23").Value = "  +2.36%  "
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.407"
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.165"
$ws.Range("E25").Value = "  -1.03%  "
$ws.Range("E27").Value = "  +1.52%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").Value = "0.0₃0737"
$ws.Range("E29").Value = "  +1.83%  "
$ws.Range("E30").Value = "  +8.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.82"
$ws.Range("E31").Value = "  -2.19%  "
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.56"
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "37.10"
$ws.Range("E35").Value = "  +2.22%  "
$ws.Range("E36").Value = "  -1.36%  "
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.835"
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.814"
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.56"
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "274.25"
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.596"
$ws.Range("E43").Value = "  +1.54%  "
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("E46").Value = "  +0.46%  "
$ws.Range("E47").Value = "  +1.34%  "
$ws.Range("D48").Value = "1.938.57"
$ws.Range("E48").Value = "  -2.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.46"
$ws.Range("E49").Value = "  +2.59%  "
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.18"
$ws.Range("E51").Value = "  -1.89%  "
